$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout:
#   A2 = "데터베스"      B2 = "데이터베이스"
# Target layout:
#   A2 = "터베"          B2 = "데이터베이스"
#   A3 = "테 터베 스"    B3 = "데이터베이스"

$ws.Range("B2").Value = "데이터베이스"
$ws.Range("A3").Value = "테 터베 스"
$ws.Range("A2").Value = "터베"
$ws.Range("B3").Value = "데이터베이스"

$ws.Range("A2").Select()
